$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 834
$ws.Range("I15").Value = 834
$ws.Range("K15").Value = 2502
$ws.Range("M15").Value = -2333

$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H80").Value = 1508.5
$ws.Range("I80").Value = 1933.1666
$ws.Range("K80").Value = 5799.4998
$ws.Range("M80").Value = -4801.4998

$ws.Range("H83").Value = 1508.5
$ws.Range("I83").Value = 1933.1666
$ws.Range("K83").Value = 17398.4994
$ws.Range("M83").Value = -12406.4994

$ws.Range("H100").Value = 1905.0526
$ws.Range("I100").Value = 1637.7858
$ws.Range("K100").Value = 1637.7858
$ws.Range("M100").Value = -1096.7858

$ws.Range("H116").Value = 23229.309
$ws.Range("I116").Value = 4230.1665
$ws.Range("K116").Value = 4230.1665
$ws.Range("M116").Value = -788.1665000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H14").Value = 385.33334
$ws.Range("I14").Value = 385.33334
$ws.Range("K14").Value = 385.33334
$ws.Range("M14").Value = -210.33334

$ws.Range("H24").Value = 0
$ws.Range("J24").Value = 0
$ws.Range("L24").Value = 0
$ws.Range("N24").ClearContents()

$ws.Range("H64").Value = 70001
$ws.Range("J64").Value = 70001
$ws.Range("L64").Value = 70001
$ws.Range("N64").Value = -70497

$ws.Range("H67").Value = 70001
$ws.Range("J67").Value = 70001
$ws.Range("L67").Value = 70001
$ws.Range("N67").Value = -71717

$ws.Range("H80").Value = 29502.715
$ws.Range("J80").Value = 33003.8
$ws.Range("L80").Value = 33003.8
$ws.Range("N80").Value = -34999.8

$ws.Range("H83").Value = 29502.715
$ws.Range("J83").Value = 33003.8
$ws.Range("L83").Value = 99011.40000000001
$ws.Range("N83").Value = -108995.4

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H102").Value = 4166.0713
$ws.Range("I102").Value = 3302.2727
$ws.Range("K102").Value = 3302.2727
$ws.Range("M102").Value = -1680.2727

$ws.Range("H122").Value = 2979.4285
$ws.Range("I122").Value = 2892.6667
$ws.Range("J122").Value = 3500
$ws.Range("K122").Value = 8678.000100000001
$ws.Range("L122").Value = 10500
$ws.Range("M122").Value = -6228.000100000001
$ws.Range("N122").Value = -15400

$ws.Range("H132").Value = 1789230.6
$ws.Range("I132").Value = 2503778.5
$ws.Range("J132").Value = 2861
$ws.Range("K132").Value = 7511335.5
$ws.Range("L132").Value = 8583
$ws.Range("M132").Value = -7508805.5
$ws.Range("N132").Value = -13643

$ws.Range("H133").Value = 62257.855
$ws.Range("J133").Value = 62257.855
$ws.Range("L133").Value = 62257.855
$ws.Range("N133").Value = -67317.85500000001

$ws.Range("H134").Value = 64990
$ws.Range("J134").Value = 64990
$ws.Range("L134").Value = 64990
$ws.Range("N134").Value = -75130

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 2000033
$ws.Range("I7").Value = 2000033
$ws.Range("K7").Value = 2000033
$ws.Range("M7").Value = -1999920

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 541.0769
$ws.Range("J19").Value = 580
$ws.Range("L19").Value = 580
$ws.Range("N19").Value = -920

$ws.Range("H22").Value = 1445.88
$ws.Range("I22").Value = 413.08334
$ws.Range("J22").Value = 2399.2307
$ws.Range("K22").Value = 413.08334
$ws.Range("L22").Value = 2399.2307
$ws.Range("M22").Value = -63.08334000000002
$ws.Range("N22").Value = -3099.2307

$ws.Range("H24").Value = 541.0769
$ws.Range("J24").Value = 580
$ws.Range("L24").Value = 580
$ws.Range("N24").Value = -920

$ws.Range("H69").Value = 15000
$ws.Range("I69").Value = 15000
$ws.Range("K69").Value = 15000
$ws.Range("M69").Value = -14251

$ws.Range("H72").Value = 15000
$ws.Range("I72").Value = 15000
$ws.Range("K72").Value = 45000
$ws.Range("M72").Value = -41256

$ws.Range("H107").Value = 570.2857
$ws.Range("I107").Value = 559.375
$ws.Range("J107").Value = 605.2
$ws.Range("K107").Value = 559.375
$ws.Range("L107").Value = 605.2
$ws.Range("M107").Value = 1360.625
$ws.Range("N107").Value = -4445.2

$ws.Range("H132").Value = 3385
$ws.Range("I132").Value = 2998.6667
$ws.Range("K132").Value = 8996.000100000001
$ws.Range("M132").Value = -6466.000100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 62500220
$ws.Range("I11").Value = 154.2
$ws.Range("K11").Value = 462.6
$ws.Range("M11").Value = -322.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H14").Value = 83751170
$ws.Range("I14").Value = 91364270
$ws.Range("K14").Value = 91364270
$ws.Range("M14").Value = -91364102

$ws.Range("H24").Value = 130124.375
$ws.Range("J24").Value = 5856.4287
$ws.Range("L24").Value = 5856.4287
$ws.Range("N24").Value = -6202.4287

$ws.Range("H113").Value = 3532.4443
$ws.Range("I113").Value = 2549.8333
$ws.Range("K113").Value = 2549.8333
$ws.Range("M113").Value = -379.8332999999998

$ws.Range("H122").Value = 4465.3335
$ws.Range("I122").Value = 4290.0835
$ws.Range("J122").Value = 5166.3335
$ws.Range("K122").Value = 12870.2505
$ws.Range("L122").Value = 15499.0005
$ws.Range("M122").Value = -10420.2505
$ws.Range("N122").Value = -20399.0005

$ws.Range("H126").Value = 3052.6
$ws.Range("I126").Value = 2838
$ws.Range("J126").Value = 3374.5
$ws.Range("K126").Value = 8514
$ws.Range("L126").Value = 10123.5
$ws.Range("M126").Value = -6044
$ws.Range("N126").Value = -15063.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4902.8335
$ws.Range("I7").Value = 3549.2856
$ws.Range("K7").Value = 3549.2856
$ws.Range("M7").Value = -3437.2856

$ws.Range("H21").Value = 1900
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 1900
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 1900
$ws.Range("N21").Value = -2248
$ws.Range("M21").ClearContents()

$ws.Range("H109").Value = 0
$ws.Range("J109").Value = 0
$ws.Range("L109").Value = 0
$ws.Range("N109").ClearContents()

$ws.Range("H126").Value = 4902.8335
$ws.Range("I126").Value = 3549.2856
$ws.Range("K126").Value = 10647.8568
$ws.Range("M126").Value = -8177.856800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I14").Value = 2122.5557
$ws.Range("J14").Value = 3540
$ws.Range("K14").Value = 2122.5557
$ws.Range("L14").Value = 3540
$ws.Range("M14").Value = -1954.5557
$ws.Range("N14").Value = -3876

$ws.Range("H22").Value = 12000
$ws.Range("J22").Value = 12000
$ws.Range("L22").Value = 12000
$ws.Range("N22").Value = -12586

$ws.Range("H32").Value = 15000
$ws.Range("I32").Value = 15000
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 15000
$ws.Range("L32").Value = 0
$ws.Range("M32").Value = -14683
$ws.Range("N32").ClearContents()

$ws.Range("H126").Value = 2521.7144
$ws.Range("I126").Value = 2830.6
$ws.Range("J126").Value = 1749.5
$ws.Range("K126").Value = 8491.799999999999
$ws.Range("L126").Value = 5248.5
$ws.Range("M126").Value = -6021.799999999999
$ws.Range("N126").Value = -10188.5
